# Rename transcript speaker labels in column D (Speaker) of the DataSheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of old speaker labels to their new abbreviated labels.
$map = @{
    "RBD"       = "T"
    "Student"   = "S"
    "Student 2" = "SN"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)  # Column D
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey([string]$val)) {
        $cell.Value = $map[[string]$val]
    }
}
